$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (cell reference, new text value) taken from the source diff.
# Column D ("Price") values look numeric (e.g. "352.25", "0.999", "3.30") but
# must stay stored as literal text, matching the original workbook, so for
# those cells we briefly force a Text number format while writing the value
# and then restore the default "Normal" style so no visible formatting changes.
$updates = @(
    @('D2', '51.514.11'),
    @('E2', '  -0.41%  '),
    @('D3', '2.810.84'),
    @('E3', '  +1.09%  '),
    @('E4', '  +0.14%  '),
    @('D5', '352.25'),
    @('E5', '  +5.57%  '),
    @('D6', '112.85'),
    @('E6', '  -3.21%  '),
    @('D7', '0.567'),
    @('E7', '  +5.05%  '),
    @('D8', '0.999'),
    @('E8', '  +0.08%  '),
    @('D9', '0.594'),
    @('E9', '  +3.05%  '),
    @('D10', '41.31'),
    @('E10', '  -1.72%  '),
    @('D11', '0.0851'),
    @('E11', '  -1.59%  '),
    @('E12', '  +1.03%  '),
    @('D13', '19.92'),
    @('E13', '  -2.22%  '),
    @('E14', '  +0.24%  '),
    @('D15', '3.260.63'),
    @('E15', '  +1.67%  '),
    @('D16', '2.816.04'),
    @('E16', '  +1.48%  '),
    @('D17', '0.885'),
    @('E17', '  -0.87%  '),
    @('D18', '51.371.99'),
    @('E18', '  -0.59%  '),
    @('D19', '7.38'),
    @('E19', '  +7.27%  '),
    @('E20', '  -4.56%  '),
    @('D21', '13.35'),
    @('E21', '  -1.79%  '),
    @('E22', '  +1.13%  '),
    @('D23', '270.50'),
    @('E23', '  -2.92%  '),
    @('D24', '69.50'),
    @('E24', '  -0.54%  '),
    @('D25', '2.74'),
    @('E25', '  +1.93%  '),
    @('D26', '26.61'),
    @('E26', '  -1.19%  '),
    @('E27', '  +0.05%  '),
    @('D28', '10.29'),
    @('E28', '  +0.81%  '),
    @('E29', '  +0.75%  '),
    @('E30', '  -2.53%  '),
    @('D32', '33.79'),
    @('E32', '  -4.30%  '),
    @('D33', '5.81'),
    @('E33', '  +3.78%  '),
    @('D34', '0.0443'),
    @('E34', '  +24.53%  '),
    @('D35', '0.0820'),
    @('E35', '  -0.34%  '),
    @('E36', '  +0.07%  '),
    @('E37', '  -2.31%  '),
    @('E38', '  -1.96%  '),
    @('E39', '  -2.36%  '),
    @('D40', '17.99'),
    @('E40', '  -5.91%  '),
    @('D41', '23.84'),
    @('E41', '  +1.89%  '),
    @('D42', '0.116'),
    @('E42', '  +1.79%  '),
    @('D43', '126.71'),
    @('E43', '  -0.52%  '),
    @('E44', '  +1.54%  '),
    @('E45', '  -0.63%  '),
    @('D46', '2.076.98'),
    @('E46', '  -0.98%  '),
    @('D47', '3.30'),
    @('E47', '  -0.86%  '),
    @('E48', '  +3.42%  '),
    @('D49', '5.64'),
    @('E49', '  +1.14%  '),
    @('E50', '  +3.85%  '),
    @('D51', '60.55'),
    @('E51', '  +0.17%  '),
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $range = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D")) {
        $range.NumberFormat = "@"
        $range.Value = $newValue
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}

Write-Output "Applied $($updates.Count) cell updates"
